$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 249324
$ws.Range("E3").Value = 1036461686

$ws.Range("C6").Value = 20791
$ws.Range("E6").Value = 360546521

$ws.Range("C7").Value = 7008
$ws.Range("E7").Value = 290033265

$ws.Range("C63").Value = 14345
$ws.Range("E63").Value = 36146424

$ws.Range("C64").Value = 5196
$ws.Range("E64").Value = 20331329

$ws.Range("C70").Value = 15723
$ws.Range("E70").Value = 24657676

$ws.Range("C74").Value = 939
$ws.Range("E74").Value = 4175586

$ws.Range("C79").Value = 116591
$ws.Range("E79").Value = 447355948

$ws.Range("C91").Value = 151092
$ws.Range("E91").Value = 482061222

$ws.Range("C92").Value = 408987
$ws.Range("D92").Value = 70904
$ws.Range("E92").Value = 1593473857

$ws.Range("C94").Value = 94135
$ws.Range("E94").Value = 915375566

$ws.Range("C95").Value = 50704
$ws.Range("E95").Value = 929563081

$ws.Range("C96").Value = 17226
$ws.Range("E96").Value = 787987719

$ws.Range("C98").Value = 808
$ws.Range("E98").Value = 117608252

$ws.Range("C102").Value = 107
$ws.Range("E102").Value = 19689236

$ws.Range("C105").Value = 8168
$ws.Range("E105").Value = 16867318

$ws.Range("C106").Value = 18336
$ws.Range("E106").Value = 41281761

$ws.Range("C107").Value = 6386
$ws.Range("E107").Value = 21935846

$ws.Range("C110").Value = 394
$ws.Range("E110").Value = 16566525

$ws.Range("C111").Value = 115
$ws.Range("E111").Value = 7718509

$ws.Range("C115").Value = 11689
$ws.Range("E115").Value = 32938514

$ws.Range("C142").Value = 168971
$ws.Range("E142").Value = 681755958

$ws.Range("C143").Value = 64956
$ws.Range("E143").Value = 373529758

$ws.Range("C145").Value = 11831
$ws.Range("E145").Value = 182701622

$ws.Range("C165").Value = 83802
$ws.Range("E165").Value = 354971515

$ws.Range("C166").Value = 35929
$ws.Range("E166").Value = 210575326

$ws.Range("C172").Value = 22700
$ws.Range("E172").Value = 44672707

$ws.Range("C175").Value = 80778
$ws.Range("E175").Value = 486144839

$ws.Range("C177").Value = 14707
$ws.Range("E177").Value = 251104307

$ws.Range("C178").Value = 4922
$ws.Range("E178").Value = 212940437
